# Rewrite the four M2Doc "complex field" placeholders
# ({m: ...}) into plain text runs containing the literal
# "{m: ...}" token text, removing the fldChar begin/instrText/fldChar
# end structure, per the TokenIteratorFieldRewriterSplit change.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Paragraph 2: " m: 2.myTemplate() " -> "{m: 2.myTemplate()}" ---
$p = $d.Paragraphs.Item(2)
$xml = '<w:p ' + $wNs + ' w:rsidP="00F5495F" w:rsidR="00735354" w:rsidRDefault="00735354" w:rsidRPr="00DC5685">' + `
  '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">{m: 2.myTemplate()}</w:t></w:r>' + `
  '</w:p>'
$p.Range.InsertXML($xml)

# --- Paragraph 4: " m:template myTemplate(a:" / ") " (split by a
#     _GoBack bookmark) -> "{m:template myTemplate(a:" / ")}" ---
$p = $d.Paragraphs.Item(4)
$xml = '<w:p ' + $wNs + ' w:rsidP="00735354" w:rsidR="00735354" w:rsidRDefault="00735354">' + `
  '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>{m:template myTemplate(a:</w:t></w:r>' + `
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">)}</w:t></w:r>' + `
  '</w:p>'
$p.Range.InsertXML($xml)

# --- Paragraph 5: " m: a + a " -> "{m: a + a}" ---
$p = $d.Paragraphs.Item(5)
$xml = '<w:p ' + $wNs + ' w:rsidP="00735354" w:rsidR="00735354" w:rsidRDefault="00735354">' + `
  '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">{m: a + a}</w:t></w:r>' + `
  '</w:p>'
$p.Range.InsertXML($xml)

# --- Paragraph 6: " m:endtemplate " -> "{m:endtemplate}" ---
$p = $d.Paragraphs.Item(6)
$xml = '<w:p ' + $wNs + ' w:rsidR="007A2DC4" w:rsidRDefault="00735354" w:rsidRPr="00DC5685">' + `
  '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">{m:endtemplate}</w:t></w:r>' + `
  '</w:p>'
$p.Range.InsertXML($xml)
